# Apply updated cryptocurrency price/volume data (and the PEPE/ImmutableX row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.510.83"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.985.02"
$ws.Range("E3").Value = "  +5.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.91"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4645"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3945"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.32"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07931"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.46"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.985.85"
$ws.Range("E13").Value = "  +9.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.196"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.855"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07118"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.83"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009960"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.17"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9973"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.547.76"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.552"
$ws.Range("E23").Value = "  +4.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.231.89"
$ws.Range("E25").Value = "  +8.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.119"
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.52"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.63"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.006"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.43"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.917"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09415"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8940"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.000004209"
$ws.Range("E34").Value = "  +153.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.281"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.347"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.168"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05818"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.179"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02129"
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.965"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9990"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5768"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1821"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.821"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5384"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.179"
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.631"
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06980"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.870"
$ws.Range("E51").Value = "  +0.65%  "
